# WAT api scripts automation
# - Renames the "typeahead/name" API path test rows to "typeahead/lastName"
# - Marks the three ORCID/RID search-cluster rows (29-31) and the two
#   invalid-ORCID/RID rows (32-33) as PASS in the STATUS column
# - Scrolls the sheet view so C11 is the top-left visible cell

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WoS_AuthorTransformation")

# The two "typeahead" rows that used to hit /typeahead/name now call the
# dedicated last-name endpoint.
$ws.Range("D27").Value = "/recommend/search/author/typeahead/lastName"
$ws.Range("D28").Value = "/recommend/search/author/typeahead/lastName"

# Rows 29-33 (ORCID / RID search-cluster tests) were missing their STATUS;
# fill them in with PASS, matching the rest of the sheet (default style,
# same as the other STATUS cells above them).
foreach ($r in 29..33) {
    $cell = $ws.Range("L" + $r)
    $cell.Style = "Normal"
    $cell.Value = "PASS"
}

# Scroll the view so the visible top-left cell is C11 (selection stays D29).
$win = $wb.Windows.Item(1)
$win.ScrollRow = 11
$win.ScrollColumn = 3
